$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": two new rows (4 & 5), one per new file handed off.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-33-13 00:33:49"

$wsOverview.Range("A5").Value = "ed5c1042-a7f0-40f7-b12a-3e96891ca660.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-33-13 00:33:49"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/e2e/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md", "", "", "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/e2e/ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "", "", "ed5c1042-a7f0-40f7-b12a-3e96891ca660.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": two new rows (4 & 5), mirroring row 3's layout.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.f3782022c08eb8a412fe4b999d71348b5081669f.zh-cn.xlf"
$wsZhCn.Range("E4").Value = "2016-03-13 00:33:45"
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I4").Value = "Include"

$wsZhCn.Range("A5").Value = "ed5c1042-a7f0-40f7-b12a-3e96891ca660.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.zh-cn.xlf"
$wsZhCn.Range("E5").Value = "2016-03-13 00:33:45"
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I5").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/e2e/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md", "", "", "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/e2e/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.f3782022c08eb8a412fe4b999d71348b5081669f.zh-cn.xlf", "", "", "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.f3782022c08eb8a412fe4b999d71348b5081669f.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/e2e/ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "", "", "ed5c1042-a7f0-40f7-b12a-3e96891ca660.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/e2e/ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.zh-cn.xlf", "", "", "ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": two new rows (4 & 5), mirroring row 3's layout.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.f3782022c08eb8a412fe4b999d71348b5081669f.de-de.xlf"
$wsDeDe.Range("E4").Value = "2016-03-13 00:33:49"
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I4").Value = "Include"

$wsDeDe.Range("A5").Value = "ed5c1042-a7f0-40f7-b12a-3e96891ca660.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.de-de.xlf"
$wsDeDe.Range("E5").Value = "2016-03-13 00:33:49"
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I5").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/e2e/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md", "", "", "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/e2e/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.f3782022c08eb8a412fe4b999d71348b5081669f.de-de.xlf", "", "", "5cb73a78-8aa3-4ac2-b47b-adaabfccc25a.f3782022c08eb8a412fe4b999d71348b5081669f.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/e2e/ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "", "", "ed5c1042-a7f0-40f7-b12a-3e96891ca660.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/e2e/ed5c1042-a7f0-40f7-b12a-3e96891ca660.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed5c1042-a7f0-40f7-b12a-3e96891ca660/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.de-de.xlf", "", "", "ed5c1042-a7f0-40f7-b12a-3e96891ca660.5ccc0c920c3098bcb5aacf10933cc96144a8f031.de-de.xlf") | Out-Null

Write-Output "Report generated for handback."
